$d = $word.ActiveDocument

# Locate the "Platform impact" bullet paragraph under KEY ACHIEVEMENTS AND IMPACT
$targetIdx = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text -like "*Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations*") {
        $targetIdx = $idx
        break
    }
}

if ($targetIdx -eq -1) {
    throw "Could not find anchor paragraph 'Platform impact...'"
}

$bullets = @(
    "• Real-time collaboration at national scale",
    "• Predictive excellence: Utilized advanced sampling methods to decrease survey margin of error from ±4.2% to ±2.1%",
    "• Increased voter turnout prediction accuracy from 71% to 87%",
    "• Executive authority: Briefed Presidents, Congressmen, Senators, Governors on election integrity, voter sentiment and postmortem analysis"
)

$insertAfterIdx = $targetIdx
foreach ($bulletText in $bullets) {
    $anchor = $d.Paragraphs.Item($insertAfterIdx)
    $anchor.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($insertAfterIdx + 1)
    $newPara.Range.Text = $bulletText
    $insertAfterIdx = $insertAfterIdx + 1
}

Write-Output "Inserted $($bullets.Count) bullets after paragraph $targetIdx"
